# Auto-applies the cryptos.xlsx data refresh described by the commit diff.
# Most cells are plain text (inline strings) holding price/percentage data;
# Range.Value on a scratch cell pre-formatted as Text (@) is used for any
# new value that would otherwise be auto-coerced to a number by Excel, so
# the destination cell's type/style stays untouched (matches original 'General').

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Scratch cell used to stage text-looking-like-numbers so they paste as text
# without ever letting Excels type inference turn them into real numbers.
$scratch = $ws.Range("Z1")
$scratch.NumberFormat = "@"

function Set-TextValue($cellAddr, $text) {
    $scratch.Value = $text
    $scratch.Copy()
    $ws.Range($cellAddr).PasteSpecial(-4163)
}

$ws.Range("D2").Value = '64.695.81'
$ws.Range("E2").Value = '  -0.23%  '
$ws.Range("D3").Value = '3.471.95'
$ws.Range("E3").Value = '  +0.41%  '
$ws.Range("E4").Value = '  +0.02%  '
Set-TextValue "D5" '575.53'
$ws.Range("E5").Value = '  -0.13%  '
Set-TextValue "D6" '160.53'
$ws.Range("E6").Value = '  -0.24%  '
Set-TextValue "D7" '0.999'
$ws.Range("E7").Value = '  +0.04%  '
$ws.Range("D8").Value = '3.476.38'
$ws.Range("E8").Value = '  +0.38%  '
$ws.Range("E9").Value = '  -6.30%  '
Set-TextValue "D10" '7.22'
$ws.Range("E10").Value = '  -0.31%  '
$ws.Range("E11").Value = '  -2.43%  '
$ws.Range("E12").Value = '  -3.19%  '
$ws.Range("D13").Value = '4.080.23'
$ws.Range("E13").Value = '  +0.81%  '
$ws.Range("E14").Value = '  -0.33%  '
$ws.Range("E15").Value = '  -2.87%  '
Set-TextValue "D16" '0.0000177'
$ws.Range("E16").Value = '  -7.70%  '
$ws.Range("B17").Value = 'WrappedEther'
$ws.Range("C17").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D17").Value = '3.542.19'
$ws.Range("E17").Value = '  +2.03%  '
$ws.Range("B18").Value = 'WrappedBTC'
$ws.Range("C18").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D18").Value = '64.756.31'
$ws.Range("E18").Value = '  -0.13%  '
Set-TextValue "D19" '6.21'
$ws.Range("E19").Value = '  -4.37%  '
Set-TextValue "D20" '13.83'
$ws.Range("E20").Value = '  -3.79%  '
Set-TextValue "D21" '381.29'
$ws.Range("E21").Value = '  -0.06%  '
$ws.Range("E22").Value = '  -1.91%  '
Set-TextValue "D23" '72.76'
$ws.Range("E23").Value = '  -0.03%  '
$ws.Range("E24").Value = '  +0.34%  '
$ws.Range("E25").Value = '  -4.00%  '
$ws.Range("E26").Value = '  +0.84%  '
$ws.Range("E27").Value = '  -1.68%  '
$ws.Range("E28").Value = '  +0.58%  '
$ws.Range("E29").Value = '  +0.31%  '
$ws.Range("B30").Value = 'NEARProtocol'
$ws.Range("C30").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextValue "D30" '6.13'
$ws.Range("E30").Value = '  -1.40%  '
$ws.Range("B31").Value = 'Fetch.AI'
$ws.Range("C31").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
Set-TextValue "D31" '1.44'
$ws.Range("E31").Value = '  -4.15%  '
$ws.Range("E32").Value = '  -1.49%  '
Set-TextValue "D33" '23.38'
$ws.Range("E33").Value = '  -0.80%  '
Set-TextValue "D34" '7.07'
$ws.Range("E34").Value = '  -2.70%  '
$ws.Range("E35").Value = '  -2.33%  '
$ws.Range("E36").Value = '  +0.04%  '
Set-TextValue "D37" '1.87'
$ws.Range("E37").Value = '  -2.63%  '
$ws.Range("E38").Value = '  -3.49%  '
Set-TextValue "D39" '26.81'
$ws.Range("E39").Value = '  +0.51%  '
$ws.Range("B40").Value = 'Mantle'
$ws.Range("C40").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
Set-TextValue "D40" '0.816'
$ws.Range("E40").Value = '  +4.87%  '
$ws.Range("B41").Value = 'Maker'
$ws.Range("C41").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D41").Value = '2.853.28'
$ws.Range("E41").Value = '  -2.87%  '
$ws.Range("E42").Value = '  -3.06%  '
$ws.Range("E43").Value = '  +0.23%  '
Set-TextValue "D44" '6.50'
$ws.Range("E44").Value = '  -3.84%  '
Set-TextValue "D45" '25.82'
$ws.Range("E45").Value = '  -0.73%  '
$ws.Range("E46").Value = '  -3.93%  '
Set-TextValue "D47" '2.43'
$ws.Range("E47").Value = '  +11.71%  '
Set-TextValue "D48" '330.52'
$ws.Range("E48").Value = '  +4.10%  '
$ws.Range("E49").Value = '  -2.42%  '
$ws.Range("B50").Value = 'SuiNetwork'
$ws.Range("C50").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
Set-TextValue "D50" '0.847'
$ws.Range("E50").Value = '  -3.66%  '
$ws.Range("B51").Value = 'Cosmos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-TextValue "D51" '6.46'
$ws.Range("E51").Value = '  -2.35%  '

# Clean up the scratch cell entirely (contents + formatting) so it leaves
# no trace and the sheet dimension stays A1:E51.
$scratch.Clear()
